$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) cells: force Text number format first so numeric-looking
# strings (e.g. "487.18", "0.999") are kept as text, matching the source data
# which preserves exact decimal formatting (trailing zeros, multi-dot big-number
# separators, etc.) instead of being reinterpreted as real numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "56.181.28"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "2.469.87"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "487.18"
$ws.Range("E5").Value = "  +5.07%  "
$ws.Range("D6").Value = "146.28"
$ws.Range("E6").Value = "  +11.45%  "
$ws.Range("D8").Value = "0.508"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").Value = "2.479.85"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "5.78"
$ws.Range("E10").Value = "  +8.94%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  +5.79%  "
$ws.Range("D14").Value = "2.911.79"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "56.224.81"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "21.00"
$ws.Range("E16").Value = "  +7.02%  "
$ws.Range("E17").Value = "  +2.62%  "
$ws.Range("D18").Value = "2.480.22"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "4.50"
$ws.Range("E19").Value = "  +8.09%  "
$ws.Range("D20").Value = "10.00"
$ws.Range("E20").Value = "  +6.67%  "
$ws.Range("D21").Value = "316.73"
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +8.65%  "
$ws.Range("D24").Value = "58.36"
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  +7.30%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  +4.55%  "
$ws.Range("D28").Value = "2.583.00"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  +8.27%  "
$ws.Range("D30").Value = "0.0₃0786"
$ws.Range("E30").Value = "  +10.67%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "149.07"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "18.15"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  +5.26%  "
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").Value = "  +8.01%  "
$ws.Range("E37").Value = "  +5.91%  "
$ws.Range("D38").Value = "0.857"
$ws.Range("E38").Value = "  +7.29%  "
$ws.Range("D39").Value = "34.16"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("D40").Value = "3.52"
$ws.Range("D41").Value = "0.996"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "0.0554"
$ws.Range("E42").Value = "  +6.72%  "
$ws.Range("D43").Value = "0.603"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  +7.53%  "
$ws.Range("D45").Value = "4.79"
$ws.Range("E45").Value = "  +13.83%  "
$ws.Range("D46").Value = "0.0924"
$ws.Range("E46").Value = "  +6.15%  "
$ws.Range("D47").Value = "257.31"
$ws.Range("E47").Value = "  +13.40%  "
$ws.Range("D48").Value = "10.17"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "0.0227"
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("D50").Value = "17.52"
$ws.Range("E50").Value = "  +6.24%  "
$ws.Range("D51").Value = "1.858.39"
$ws.Range("E51").Value = "  -3.80%  "
